## Add files via upload
## Adds a new slide 6 ("Naslov 1" / title-only layout) with the
## bold+italic title "Preizkus aplikacije", positioned/sized to match
## the authored slide, as the new last slide of the deck.

$p = $ppt.ActivePresentation

# EMU -> points; PowerPoint shape geometry (Left/Top/Width/Height) is
# stored internally as a single-precision Point value that gets
# multiplied by 12700 to obtain EMU on save. Nudge by half an EMU
# (in point-space) before the float32 truncation so the EMU value
# that gets serialized lands exactly on the target integer instead of
# rounding down.
function EmuToPt($emu) {
    return ($emu + 0.5) / 12700
}

# ppLayoutTitleOnly = 11 -> slideLayout6.xml ("Samo naslov" / Title Only),
# the same "title placeholder only" layout used by the authored slide.
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 11)

$sh = $s.Shapes.Item(1)
$sh.Name = "Naslov 1"

$sh.Left   = EmuToPt 3656604
$sh.Top    = EmuToPt 3034529
$sh.Width  = EmuToPt 5150511
$sh.Height = EmuToPt 788941

$tr = $sh.TextFrame.TextRange
$tr.Text = "Preizkus aplikacije"
$tr.Font.Bold = $true
$tr.Font.Italic = $true
$tr.LanguageID = "sl-SI"
